$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(261).Resize(2).Insert()

$ws.Range("A261:A262").Value = 7
$ws.Range("B261:B262").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C261:C262").Value = "Ñuble"
$ws.Range("D261").Value = 44858
$ws.Range("D262").Value = 44858
$ws.Range("E261:E262").Value = 16
$ws.Range("F261:F262").Value = 100112006
$ws.Range("G261:G262").Value = "Repollo"
$ws.Range("H261:H262").Value = "Crespo record"
$ws.Range("I261").Value = "Primera"
$ws.Range("I262").Value = "Segunda"
$ws.Range("J261").Value = 400
$ws.Range("J262").Value = 300
$ws.Range("K261").Value = 1500
$ws.Range("K262").Value = 1200
$ws.Range("L261").Value = 1600
$ws.Range("L262").Value = 1200
$ws.Range("M261").Value = 1550
$ws.Range("M262").Value = 1200
$ws.Range("N261:N262").Value = "$/unidad"
$ws.Range("O261:O262").Value = "Provincia de Diguillín"
$ws.Range("P261").Value = 1550
$ws.Range("P262").Value = 1200
$ws.Range("Q261:Q262").Value = 1
$ws.Range("R261:R262").Value = "Hortaliza"
